$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right after "Preduslovi" section;
# it needs to move to just after the newly inserted sentence fragment.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Range.Delete()

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "gosti nemaju pravo prijavljivanja grešaka"
$found = $find.Execute()

if ($found) {
    $insertionPoint = $find.Parent
    $insertionPoint.Collapse(0)
    $insertionPoint.InsertAfter(" i drugih korisnika")

    # Re-drop the _GoBack bookmark immediately after the inserted text,
    # matching where Word leaves it after the last edit.
    $insertionPoint.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $insertionPoint) | Out-Null
}
